$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Date-row => open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
# Column I (fixed_ticker) is set to "HUBB" for every data row (2-43).
$rowData = @(
    @{Row=2; D=87.20302093263906; E=86.80421447753906; F=90.0026378660225; G=85.95874065472063; H=53140009}
    @{Row=3; D=87.16911948347052; E=83.70576477050781; F=87.7062587149555; G=81.89391532542628; H=53140009}
    @{Row=4; D=68.4535831246904; E=78.09788513183594; F=78.5817117983248; G=67.61494736508776; H=53140009}
    @{Row=5; D=82.00073373183702; E=73.39000701904297; F=83.50213150406721; G=67.489917415344; H=53140009}
    @{Row=6; D=85.70805843277698; E=86.39424896240234; F=90.8626362330204; G=85.1362381850669; H=53140009}
    @{Row=7; D=86.69585892915084; E=88.61056518554688; F=88.67630762591835; G=84.94550877368353; H=53140009}
    @{Row=8; D=88.46136560664492; E=86.3948974609375; F=89.0482481096229; G=83.60930243659449; H=53140009}
    @{Row=9; D=97.06302486520777; E=101.5377426147461; F=103.9663961048092; G=95.71562443287678; H=53140009}
    @{Row=10; D=100.7757124948134; E=94.64350128173828; F=102.5492040584174; G=94.55147579637141; H=53140009}
    @{Row=11; D=95.44903346952574; E=99.97698974609376; F=102.1988871269707; G=94.44749410623454; H=53140009}
    @{Row=12; D=98.3710181345408; E=106.5700149536133; F=107.8573738435588; G=97.13439339569244; H=53140009}
    @{Row=13; D=115.3378110606258; E=115.8576583862305; F=127.0045389882988; G=113.2584217582072; H=53140009}
    @{Row=14; D=104.3027128209891; E=89.02029418945312; F=105.2626809919584; G=87.21177431684011; H=53140009}
    @{Row=15; D=90.32077537168622; E=106.3943023681641; F=107.0589940976965; G=89.9236858651348; H=53140009}
    @{Row=16; D=116.9279155180521; E=88.32777404785156; F=117.092935117036; G=84.26313273119725; H=53140009}
    @{Row=17; D=86.9587854449486; E=95.70368957519533; F=96.80664504148184; G=84.10509938372284; H=53140009}
    @{Row=18; D=105.3055159150201; E=112.4904479980469; F=113.7533232147274; G=104.0801110879951; H=53140009}
    @{Row=19; D=117.3330058352941; E=115.3262557983398; F=118.8780344082158; G=109.9275506854532; H=53140009}
    @{Row=20; D=118.3112090793775; E=126.6598510742188; F=129.0368106776354; G=110.677650381917; H=53140009}
    @{Row=21; D=132.9430805550957; E=128.8150177001953; F=134.7957554827512; G=128.6081692424863; H=53140009}
    @{Row=22; D=99.15965838792528; E=112.6592102050781; F=121.7802308504158; G=96.977642403218; H=53140009}
    @{Row=23; D=115.1247694993675; E=123.0958557128906; F=123.5473030491343; G=111.4766742308017; H=53140009}
    @{Row=24; D=126.8520953294007; E=133.5328674316406; F=141.5167697101922; G=122.5206038781264; H=53140009}
    @{Row=25; D=144.7570471018314; E=143.6583862304688; F=158.9567038821254; G=141.87649796301; H=53140009}
    @{Row=26; D=173.6322656651773; E=178.2555236816406; F=182.3960416502286; G=169.9745019102035; H=53140009}
    @{Row=27; D=175.5247510295333; E=187.0683898925781; F=187.3576782404506; G=169.449653005908; H=53140009}
    @{Row=28; D=170.0353883733823; E=186.9310913085937; F=187.6389906422147; G=166.9600321697348; H=53140009}
    @{Row=29; D=197.4651722615401; E=176.5141754150391; F=199.6611163411355; G=169.0121658129708; H=53140009}
    @{Row=30; D=175.6006678632438; E=185.2240447998047; F=193.2735593007433; G=170.3504668221285; H=53140009}
    @{Row=31; D=169.7001169878106; E=208.8308868408203; F=210.2134278996782; G=168.9754784901359; H=53140009}
    @{Row=32; D=214.3964115464429; E=227.5626068115234; F=232.5981577271976; G=201.1966597836754; H=53140009}
    @{Row=33; D=228.6368774338058; E=220.3396148681641; F=230.3117371345561; G=212.3214882422306; H=53140009}
    @{Row=34; D=235.5010282972273; E=260.3872985839844; F=260.4356283998886; G=212.4807521388731; H=53140009}
    @{Row=35; D=320.2585651155812; E=302.8174438476562; F=330.0516000497871; G=296.3340518380372; H=53140009}
    @{Row=36; D=305.152343887533; E=263.0506286621094; F=309.4569987730827; G=241.8877522399686; H=53140009}
    @{Row=37; D=321.6758561027003; E=328.16943359375; F=339.6211759405795; G=308.419801712593; H=53140009}
    @{Row=38; D=407.8078639524077; E=363.5555725097656; F=421.5348936940516; G=361.6226000848657; H=53140009}
    @{Row=39; D=366.1117601473204; E=389.4292907714844; F=397.244454596028; G=341.6721690696752; H=53140009}
    @{Row=40; D=423.5910700222979; E=421.6164245605469; F=455.916006589339; G=414.7841380053809; H=53140009}
    @{Row=41; D=414.7728551851191; E=418.8523559570313; F=472.7127109714594; G=401.1877042374882; H=53140009}
    @{Row=42; D=327.03437204524; E=360.9004516601562; F=363.1264126640976; G=297.5456124018623; H=53140009}
    @{Row=43; D=404.2594231913012; E=436.1846313476562; F=446.4591109088868; G=400.849566860909; H=53140009}
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D   # D: open_price
    $ws.Cells.Item($r, 5).Value = $item.E   # E: close_price
    $ws.Cells.Item($r, 6).Value = $item.F   # F: high_price
    $ws.Cells.Item($r, 7).Value = $item.G   # G: low_price
    $ws.Cells.Item($r, 8).Value = $item.H   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "HUBB"    # I: fixed_ticker
}

Write-Host "Applied all changes"
